# The commit renames the display "name" of the embedded logo pictures
# (the wp:docPr/name and pic:cNvPr/name attributes) in the document's
# headers/footers:
#   - the two Pearson logo pictures (in footer1.xml / footer2.xml):
#       name="image1.png"  ->  name="image2.png"
#   - the BTEC logo picture (in header1.xml):
#       name="image2.jpg"  ->  name="image1.jpg"
#
# Word's InlineShape object doesn't expose a settable "Name" property
# (only floating Shape objects do), so we go through the document's
# flat-OPC WordOpenXML representation, patch the attribute text, and
# write it back.

$d = $word.ActiveDocument

$xml = $d.WordOpenXML

$xml = $xml.Replace('name="image1.png"', 'name="image2.png"')
$xml = $xml.Replace('name="image2.jpg"', 'name="image1.jpg"')

$d.WordOpenXML = $xml
